# SquareDeal.docx edit
#
# Commit summary: the -W (Wizard) option now replaces the "owner string"
# (instead of asking for it separately), and the other half of the
# session key is now passed base64 encoded. This requires two small
# text insertions into the paragraph that documents the -W option:
#
#   1. "...with a string as argument" + " to replace the owner string"
#      (inserted right before the following ". This string is...")
#
#   2. "...session key, and the" + " (base64 encoded)"
#      (inserted right before the following "DI and the original/...")

$d = $word.ActiveDocument

# --- Insertion 1: " to replace the owner string" -------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "with a string as argument",  # Find what
    $true,                        # MatchCase
    $false,                       # MatchWholeWord
    $false,                       # MatchWildcards
    $false,                       # MatchSoundsLike
    $false,                       # MatchAllWordForms
    $true,                        # Forward
    1,                            # Wrap (wdFindContinue)
    $false,                       # Format
    "",                           # ReplaceWith
    0)                            # Replace (wdReplaceNone)

if (-not $found1) {
    throw "Could not find anchor text 'with a string as argument'"
}

$r1.Collapse(0)                               # collapse to end of match
$r1.InsertAfter(" to replace the owner string")

# --- Insertion 2: " (base64 encoded)" -------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "session key, and the ",      # Find what (includes trailing space)
    $true,                        # MatchCase
    $false,                       # MatchWholeWord
    $false,                       # MatchWildcards
    $false,                       # MatchSoundsLike
    $false,                       # MatchAllWordForms
    $true,                        # Forward
    1,                            # Wrap (wdFindContinue)
    $false,                       # Format
    "",                           # ReplaceWith
    0)                            # Replace (wdReplaceNone)

if (-not $found2) {
    throw "Could not find anchor text 'session key, and the '"
}

$r2.Collapse(0)                               # collapse to end of match
$r2.InsertAfter("(base64 encoded) ")
